$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31 ("2.14" / Timepoint) gets "marked/resolved" styling, matching the
# highlighted rows above it (e.g. row 8): copy that formatting over, then
# bump the weight for D31 from 30 to 100 (the shared formula in E31 and the
# totals in row 35 recalculate automatically).
$ws.Range("A8:E8").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)
$ws.Range("D31").Value = 100

# Move the view: drop the scrolled-down A10 top-left anchor and select B2.
$ws.Range("B2").Select()
